$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 995.51514
$ws.Range("I15").Value = 995.51514
$ws.Range("K15").Value = 2986.54542
$ws.Range("M15").Value = -2817.54542
$ws.Range("H31").Value = 2463
$ws.Range("I31").Value = 2463
$ws.Range("K31").Value = 7389
$ws.Range("M31").Value = -7159
$ws.Range("H99").Value = 2889.3635
$ws.Range("J99").Value = 4111.3335
$ws.Range("L99").Value = 12334.0005
$ws.Range("N99").Value = -15330.0005
$ws.Range("H101").Value = 1437.4286
$ws.Range("I101").Value = 1044.0834
$ws.Range("J101").Value = 3797.5
$ws.Range("K101").Value = 3132.2502
$ws.Range("L101").Value = 11392.5
$ws.Range("M101").Value = -1510.2502
$ws.Range("N101").Value = -14636.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11450.5
$ws.Range("I2").Value = 18185.572
$ws.Range("J2").Value = 4715.4287
$ws.Range("K2").Value = 18185.572
$ws.Range("L2").Value = 4715.4287
$ws.Range("M2").Value = -18072.572
$ws.Range("N2").Value = -4941.4287
$ws.Range("H110").Value = 5359.231
$ws.Range("I110").Value = 4528.778
$ws.Range("J110").Value = 7227.75
$ws.Range("K110").Value = 4528.778
$ws.Range("L110").Value = 7227.75
$ws.Range("M110").Value = -2483.778
$ws.Range("N110").Value = -11317.75
$ws.Range("H116").Value = 11450.5
$ws.Range("I116").Value = 18185.572
$ws.Range("J116").Value = 4715.4287
$ws.Range("K116").Value = 18185.572
$ws.Range("L116").Value = 4715.4287
$ws.Range("M116").Value = -15891.572
$ws.Range("N116").Value = -9303.4287

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11450.5
$ws.Range("I3").Value = 18185.572
$ws.Range("J3").Value = 4715.4287
$ws.Range("K3").Value = 18185.572
$ws.Range("L3").Value = 4715.4287
$ws.Range("M3").Value = -18071.572
$ws.Range("N3").Value = -4943.4287
$ws.Range("H99").Value = 2993
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 1574.1936
$ws.Range("I105").Value = 1595.5186
$ws.Range("J105").Value = 1430.25
$ws.Range("K105").Value = 1595.5186
$ws.Range("L105").Value = 1430.25
$ws.Range("M105").Value = 151.4813999999999
$ws.Range("N105").Value = -4924.25
$ws.Range("H107").Value = 14613.333
$ws.Range("I107").Value = 20418.334
$ws.Range("K107").Value = 20418.334
$ws.Range("M107").Value = -18498.334
$ws.Range("H134").Value = 31251862
$ws.Range("I134").Value = 41667984
$ws.Range("K134").Value = 125003952
$ws.Range("M134").Value = -125001417
$ws.Range("H135").Value = 109999.5
$ws.Range("J135").Value = 109999.5
$ws.Range("L135").Value = 109999.5
$ws.Range("N135").Value = -120139.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2514.375
$ws.Range("I31").Value = 2180.8572
$ws.Range("K31").Value = 2180.8572
$ws.Range("M31").Value = -1885.8572
$ws.Range("H33").Value = 11596.583
$ws.Range("J33").Value = 23727
$ws.Range("L33").Value = 23727
$ws.Range("N33").Value = -24485
$ws.Range("H34").Value = 2514.375
$ws.Range("I34").Value = 2180.8572
$ws.Range("K34").Value = 2180.8572
$ws.Range("M34").Value = -1978.8572
$ws.Range("H134").Value = 2849
$ws.Range("I134").Value = 2999
$ws.Range("J134").Value = 2699
$ws.Range("K134").Value = 8997
$ws.Range("L134").Value = 8097
$ws.Range("M134").Value = -6462
$ws.Range("N134").Value = -13167

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H49").Value = 2150
$ws.Range("I49").Value = 1225
$ws.Range("K49").Value = 3675
$ws.Range("M49").Value = -3519
$ws.Range("H109").Value = 4250
$ws.Range("I109").Value = 3428.5715
$ws.Range("K109").Value = 10285.7145
$ws.Range("M109").Value = -9245.7145

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1403.6154
$ws.Range("I97").Value = 1265.9524
$ws.Range("J97").Value = 1981.8
$ws.Range("K97").Value = 1265.9524
$ws.Range("L97").Value = 1981.8
$ws.Range("M97").Value = -769.9523999999999
$ws.Range("N97").Value = -2973.8
$ws.Range("H102").Value = 2330.5789
$ws.Range("I102").Value = 2159.4443
$ws.Range("K102").Value = 2159.4443
$ws.Range("M102").Value = -537.4443000000001
$ws.Range("H132").Value = 3877.4443
$ws.Range("J132").Value = 3737.125
$ws.Range("L132").Value = 11211.375
$ws.Range("N132").Value = -16271.375
$ws.Range("H136").Value = 233866.25
$ws.Range("J136").Value = 233866.25
$ws.Range("L136").Value = 701598.75
$ws.Range("N136").Value = -706698.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 10000000
$ws.Range("I42").Value = 10000000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 10000000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -9999437
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 10000000
$ws.Range("I49").Value = 10000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 10000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -9999853
$ws.Range("N49").ClearContents()
$ws.Range("H61").Value = 3887.6
$ws.Range("I61").Value = 4110
$ws.Range("K61").Value = 4110
$ws.Range("M61").Value = -3908
$ws.Range("H113").Value = 3887.6
$ws.Range("I113").Value = 4110
$ws.Range("K113").Value = 4110
$ws.Range("M113").Value = -1940
$ws.Range("H122").Value = 9799.4
$ws.Range("I122").Value = 9749.5
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 29248.5
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -26798.5
$ws.Range("N122").Value = -34897

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 26089.5
$ws.Range("J28").Value = 27307.4
$ws.Range("L28").Value = 27307.4
$ws.Range("N28").Value = -28003.4
$ws.Range("H62").Value = 4423.077
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 4423.077
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -23740
$ws.Range("H81").Value = 36166.332
$ws.Range("I81").Value = 42199.6
$ws.Range("K81").Value = 84399.2
$ws.Range("M81").Value = -83338.2
$ws.Range("H84").Value = 36166.332
$ws.Range("I84").Value = 42199.6
$ws.Range("K84").Value = 421996
$ws.Range("M84").Value = -416692
$ws.Range("H132").Value = 6650.8
$ws.Range("I132").Value = 6438.5
$ws.Range("K132").Value = 19315.5
$ws.Range("M132").Value = -16785.5
